$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the formatting (bold, border, centered)
# already used by the other header cells (copy format from G1 -> H1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# Add corresponding data value in H2 (era data updated)
$ws.Range("H2").Value = 0
